$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.895.07"
$ws.Range("D3").Value = "1.552.83"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'206.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "'0.484"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "'21.69"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.54%  "
$ws.Range("E9").Value = "  +1.61%  "
$ws.Range("D10").Value = "'0.0586"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("D12").Value = "1.775.26"
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").Value = "1.557.26"
$ws.Range("E13").Value = "  +1.81%  "
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "'61.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "26.885.47"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").Value = "'215.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").Value = "0.0₃0689"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").Value = "'7.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").Value = "'9.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").Value = "'1.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").Value = "'152.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").Value = "'6.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.07%  "
$ws.Range("D27").Value = "'14.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").Value = "1.415.89"
$ws.Range("E33").Value = "  +4.49%  "
$ws.Range("E34").Value = "  +3.29%  "
$ws.Range("D35").Value = "'1.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.45%  "
$ws.Range("D36").Value = "'0.958"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.58%  "
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("E40").Value = "  +1.34%  "
$ws.Range("D42").Value = "'5.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("D44").Value = "'2.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.01%  "
$ws.Range("D45").Value = "'63.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.33%  "
$ws.Range("E46").Value = "  +1.64%  "
$ws.Range("D47").Value = "1.688.70"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("D48").Value = "'86.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("E49").Value = "  +1.94%  "
$ws.Range("D50").Value = "'0.0960"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.70%  "
$ws.Range("D51").Value = "'1.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.29%  "
